$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Write new header cells X1:AM1 (new fitted-parameter columns)
$ws.Range("X1").Value = "he_beta1_opt"
$ws.Range("Y1").Value = "he_beta2_opt"
$ws.Range("Z1").Value = "he_beta3_opt"
$ws.Range("AA1").Value = "he_beta4_opt"
$ws.Range("AB1").Value = "he_beta1_err"
$ws.Range("AC1").Value = "he_beta2_err"
$ws.Range("AD1").Value = "he_beta3_err"
$ws.Range("AE1").Value = "he_beta4_err"
$ws.Range("AF1").Value = "ne_beta1_opt"
$ws.Range("AG1").Value = "ne_beta2_opt"
$ws.Range("AH1").Value = "ne_beta3_opt"
$ws.Range("AI1").Value = "ne_beta4_opt"
$ws.Range("AJ1").Value = "ne_beta1_err"
$ws.Range("AK1").Value = "ne_beta2_err"
$ws.Range("AL1").Value = "ne_beta3_err"
$ws.Range("AM1").Value = "ne_beta4_err"

# Apply the same header style (bold, centered, bordered) used by the
# existing header cells to the newly added header cells.
$ws.Range("W1").Copy() | Out-Null
$ws.Range("X1:AM1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Write new data values for rows 2-15 across columns X:AM
# Row 2
$ws.Range("X2").Value = -0.2818322083709183
$ws.Range("Y2").Value = 1.61283961856221
$ws.Range("Z2").Value = -0.4023975208795856
$ws.Range("AA2").Value = 0.3789880073917434
$ws.Range("AB2").Value = 0.02899332980916463
$ws.Range("AC2").Value = 0.04860475600291164
$ws.Range("AD2").Value = 0.03983576712972371
$ws.Range("AE2").Value = 0.04525204280893394
$ws.Range("AF2").Value = -0.1464541511488946
$ws.Range("AG2").Value = 0.5088777520672764
$ws.Range("AH2").Value = -0.2067280812844753
$ws.Range("AI2").Value = 0.1927191810822214
$ws.Range("AJ2").Value = 0.0111508479665389
$ws.Range("AK2").Value = 0.01509873968252844
$ws.Range("AL2").Value = 0.01530313988961163
$ws.Range("AM2").Value = 0.01745951347565698

# Row 3
$ws.Range("X3").Value = -0.424225213963524
$ws.Range("Y3").Value = 1.622952398096494
$ws.Range("Z3").Value = -0.533863982712322
$ws.Range("AA3").Value = 0.3455764253674427
$ws.Range("AB3").Value = 0.02981935954575577
$ws.Range("AC3").Value = 0.04905079769168804
$ws.Range("AD3").Value = 0.04066662138672135
$ws.Range("AE3").Value = 0.04545260045017462
$ws.Range("AF3").Value = -0.2536610928569659
$ws.Range("AG3").Value = 0.5285288241087726
$ws.Range("AH3").Value = -0.2236841670175525
$ws.Range("AI3").Value = 0.1609364453392649
$ws.Range("AJ3").Value = 0.01074173766443454
$ws.Range("AK3").Value = 0.01444279457417761
$ws.Range("AL3").Value = 0.01461464636781925
$ws.Range("AM3").Value = 0.01664241368310447

# Row 4
$ws.Range("X4").Value = -0.3209349929185176
$ws.Range("Y4").Value = 1.601079159029586
$ws.Range("Z4").Value = -0.3014408651871817
$ws.Range("AA4").Value = 0.3520096364358008
$ws.Range("AB4").Value = 0.03013242919011222
$ws.Range("AC4").Value = 0.05010838014938088
$ws.Range("AD4").Value = 0.04084921747746115
$ws.Range("AE4").Value = 0.04670679797728246
$ws.Range("AF4").Value = -0.2181754225316178
$ws.Range("AG4").Value = 0.5191097221149259
$ws.Range("AH4").Value = -0.05418722607571241
$ws.Range("AI4").Value = 0.1539088782359884
$ws.Range("AJ4").Value = 0.01126400007120872
$ws.Range("AK4").Value = 0.01518281797579301
$ws.Range("AL4").Value = 0.01529564400653448
$ws.Range("AM4").Value = 0.0175126280232572

# Row 5
$ws.Range("X5").Value = -0.06391430267323721
$ws.Range("Y5").Value = 1.58867307157299
$ws.Range("Z5").Value = 0.1068486113443717
$ws.Range("AA5").Value = 0.3472367198769921
$ws.Range("AB5").Value = 0.03327124834479004
$ws.Range("AC5").Value = 0.05637551360347785
$ws.Range("AD5").Value = 0.04566818989448151
$ws.Range("AE5").Value = 0.05269601285326225
$ws.Range("AF5").Value = -0.07704346439376698
$ws.Range("AG5").Value = 0.5080200361723906
$ws.Range("AH5").Value = 0.1706504707896339
$ws.Range("AI5").Value = 0.1627427385918033
$ws.Range("AJ5").Value = 0.01080771817093802
$ws.Range("AK5").Value = 0.01468371026015468
$ws.Range("AL5").Value = 0.01485994416889032
$ws.Range("AM5").Value = 0.0169648520882047

# Row 6
$ws.Range("X6").Value = 0.142608331342784
$ws.Range("Y6").Value = 1.556822479916655
$ws.Range("Z6").Value = 0.3532707038775365
$ws.Range("AA6").Value = 0.3548492741559568
$ws.Range("AB6").Value = 0.03492989783779617
$ws.Range("AC6").Value = 0.05850319931763003
$ws.Range("AD6").Value = 0.04841484444997164
$ws.Range("AE6").Value = 0.05514756439594747
$ws.Range("AF6").Value = 0.0585000877007194
$ws.Range("AG6").Value = 0.4934741716151153
$ws.Range("AH6").Value = 0.2652043287443192
$ws.Range("AI6").Value = 0.1825874467317988
$ws.Range("AJ6").Value = 0.01130059842508973
$ws.Range("AK6").Value = 0.01533468228849906
$ws.Range("AL6").Value = 0.01562275089409195
$ws.Range("AM6").Value = 0.01775973475815213

# Row 7
$ws.Range("X7").Value = 0.08823447564722595
$ws.Range("Y7").Value = 1.590179851659439
$ws.Range("Z7").Value = 0.1813037574670601
$ws.Range("AA7").Value = 0.3418133277088858
$ws.Range("AB7").Value = 0.03913758448539956
$ws.Range("AC7").Value = 0.06628516528276628
$ws.Range("AD7").Value = 0.05381338244388821
$ws.Range("AE7").Value = 0.06191470271589609
$ws.Range("AF7").Value = 0.06105710247052874
$ws.Range("AG7").Value = 0.501787791378859
$ws.Range("AH7").Value = 0.1278423015209272
$ws.Range("AI7").Value = 0.1876090544954163
$ws.Range("AJ7").Value = 0.01081805663271417
$ws.Range("AK7").Value = 0.01469376552970192
$ws.Range("AL7").Value = 0.01485885606616302
$ws.Range("AM7").Value = 0.01700310077525933

# Row 8
$ws.Range("X8").Value = -0.16897006386758
$ws.Range("Y8").Value = 1.59665978619713
$ws.Range("Z8").Value = -0.2133495590489426
$ws.Range("AA8").Value = 0.3288822022274895
$ws.Range("AB8").Value = 0.03024062646451737
$ws.Range("AC8").Value = 0.05106140344991802
$ws.Range("AD8").Value = 0.04144927237866335
$ws.Range("AE8").Value = 0.04758023476900073
$ws.Range("AF8").Value = -0.08054288970429216
$ws.Range("AG8").Value = 0.5218440921152593
$ws.Range("AH8").Value = -0.1132129287907175
$ws.Range("AI8").Value = 0.1603015537785992
$ws.Range("AJ8").Value = 0.01120815260102429
$ws.Range("AK8").Value = 0.01525233994765576
$ws.Range("AL8").Value = 0.01537857508312644
$ws.Range("AM8").Value = 0.01758995261205635

# Row 9
$ws.Range("X9").Value = -0.3742581481674317
$ws.Range("Y9").Value = 1.64417733494512
$ws.Range("Z9").Value = -0.4881323809835809
$ws.Range("AA9").Value = 0.2972696944276134
$ws.Range("AB9").Value = 0.03105884568035046
$ws.Range("AC9").Value = 0.05179721098981729
$ws.Range("AD9").Value = 0.0424899052895765
$ws.Range("AE9").Value = 0.04760846967449125
$ws.Range("AF9").Value = -0.2088132090097514
$ws.Range("AG9").Value = 0.5417999249063767
$ws.Range("AH9").Value = -0.2114331559701047
$ws.Range("AI9").Value = 0.1378996784138836
$ws.Range("AJ9").Value = 0.01034555246760838
$ws.Range("AK9").Value = 0.01399810970295425
$ws.Range("AL9").Value = 0.01413130185952091
$ws.Range("AM9").Value = 0.01609204730473222

# Row 10
$ws.Range("X10").Value = -0.3741845074124353
$ws.Range("Y10").Value = 1.645986454292909
$ws.Range("Z10").Value = -0.4179257166716254
$ws.Range("AA10").Value = 0.2877239971599133
$ws.Range("AB10").Value = 0.0299757341879162
$ws.Range("AC10").Value = 0.05001424225861573
$ws.Range("AD10").Value = 0.04071093286955629
$ws.Range("AE10").Value = 0.04592551846498218
$ws.Range("AF10").Value = -0.2357687091041426
$ws.Range("AG10").Value = 0.550451902551256
$ws.Range("AH10").Value = -0.1343464149348503
$ws.Range("AI10").Value = 0.1282080444387994
$ws.Range("AJ10").Value = 0.009872744581635652
$ws.Range("AK10").Value = 0.01333824131178331
$ws.Range("AL10").Value = 0.01340704549156412
$ws.Range("AM10").Value = 0.01531231222443566

# Row 11
$ws.Range("X11").Value = -0.1795981181602798
$ws.Range("Y11").Value = 1.613472805312679
$ws.Range("Z11").Value = -0.07340936480770488
$ws.Range("AA11").Value = 0.2959324905326929
$ws.Range("AB11").Value = 0.03050986031923004
$ws.Range("AC11").Value = 0.05169066957304505
$ws.Range("AD11").Value = 0.04158384927975559
$ws.Range("AE11").Value = 0.04787495991505471
$ws.Range("AF11").Value = -0.1375631923556189
$ws.Range("AG11").Value = 0.5358542651341202
$ws.Range("AH11").Value = 0.07277856451303014
$ws.Range("AI11").Value = 0.1296782192522881
$ws.Range("AJ11").Value = 0.01046570806201321
$ws.Range("AK11").Value = 0.01422782450250441
$ws.Range("AL11").Value = 0.01430683462342983
$ws.Range("AM11").Value = 0.0163655401313579

# Row 12
$ws.Range("X12").Value = 0.05911902069075572
$ws.Range("Y12").Value = 1.610383385803687
$ws.Range("Z12").Value = 0.2510136334248523
$ws.Range("AA12").Value = 0.2862691824062578
$ws.Range("AB12").Value = 0.03701042611911494
$ws.Range("AC12").Value = 0.06306141084308856
$ws.Range("AD12").Value = 0.05112067178388197
$ws.Range("AE12").Value = 0.058421792337204
$ws.Range("AF12").Value = 0.0003536618119748749
$ws.Range("AG12").Value = 0.5252390579108185
$ws.Range("AH12").Value = 0.232167865028879
$ws.Range("AI12").Value = 0.1413206552478206
$ws.Range("AJ12").Value = 0.01117619904171851
$ws.Range("AK12").Value = 0.01523775393883246
$ws.Range("AL12").Value = 0.01543264990907845
$ws.Range("AM12").Value = 0.01755630925589435

# Row 13
$ws.Range("X13").Value = 0.1346780270657072
$ws.Range("Y13").Value = 1.607185313339139
$ws.Range("Z13").Value = 0.3024361371652047
$ws.Range("AA13").Value = 0.3125292812530701
$ws.Range("AB13").Value = 0.03332679473004514
$ws.Range("AC13").Value = 0.05655263956172073
$ws.Range("AD13").Value = 0.04603548543429755
$ws.Range("AE13").Value = 0.05250861214328094
$ws.Range("AF13").Value = 0.0689821561900649
$ws.Range("AG13").Value = 0.5192875359947481
$ws.Range("AH13").Value = 0.2119816911058096
$ws.Range("AI13").Value = 0.1504557300873637
$ws.Range("AJ13").Value = 0.0102735673446443
$ws.Range("AK13").Value = 0.01398145233701091
$ws.Range("AL13").Value = 0.01415592644342237
$ws.Range("AM13").Value = 0.01612499614996098

# Row 14
$ws.Range("X14").Value = -0.05347759461339245
$ws.Range("Y14").Value = 1.598020345393183
$ws.Range("Z14").Value = -0.04176900819463705
$ws.Range("AA14").Value = 0.2758947538861392
$ws.Range("AB14").Value = 0.03695448740835623
$ws.Range("AC14").Value = 0.06278042821260278
$ws.Range("AD14").Value = 0.05067898764801088
$ws.Range("AE14").Value = 0.05831100637097937
$ws.Range("AF14").Value = -0.01998885829008364
$ws.Range("AG14").Value = 0.5311886394751772
$ws.Range("AH14").Value = -0.01627167409435862
$ws.Range("AI14").Value = 0.1466547519352715
$ws.Range("AJ14").Value = 0.0116070392657936
$ws.Range("AK14").Value = 0.01583584964972957
$ws.Range("AL14").Value = 0.01592379408625677
$ws.Range("AM14").Value = 0.01823403306086632

# Row 15
$ws.Range("X15").Value = -0.2962098041367572
$ws.Range("Y15").Value = 1.641073670275106
$ws.Range("Z15").Value = -0.3871556532407079
$ws.Range("AA15").Value = 0.2664300714454093
$ws.Range("AB15").Value = 0.03231309941190843
$ws.Range("AC15").Value = 0.0544567507805716
$ws.Range("AD15").Value = 0.04425486433251795
$ws.Range("AE15").Value = 0.05001309943663997
$ws.Range("AF15").Value = -0.1606111143257737
$ws.Range("AG15").Value = 0.5484845488378838
$ws.Range("AH15").Value = -0.1834289342750187
$ws.Range("AI15").Value = 0.1328434387706706
$ws.Range("AJ15").Value = 0.0107133235112358
$ws.Range("AK15").Value = 0.01456620343432114
$ws.Range("AL15").Value = 0.01467248999036309
$ws.Range("AM15").Value = 0.01672824529988068

